# Update the "Förändrad" (last-changed) date in column C for every data row,
# and append the record's friendly name as a second HYPERLINK() argument for
# the rows that already carry HYPERLINK formulas in columns S, T, U, V, W, X, Y.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 319

for ($r = $firstRow; $r -le $lastRow; $r++) {
    # Bump the "changed" timestamp column (C) for every row.
    $ws.Cells.Item($r, 3).Value = 45186

    # Column A holds the record id (e.g. "A 44561-2021") used as the
    # friendly hyperlink text.
    $idCell = $ws.Cells.Item($r, 1)
    $id = $idCell.Value2

    foreach ($col in 19, 20, 21, 22, 23, 24, 25) {
        $cell = $ws.Cells.Item($r, $col)
        $formula = $cell.Formula
        if ($formula -and $formula.Length -gt 0 -and $formula.ToUpper().StartsWith("=HYPERLINK(")) {
            if ($formula.IndexOf(",") -lt 0) {
                $newFormula = $formula.Substring(0, $formula.Length - 1) + ', "' + $id + '")'
                $cell.Formula = $newFormula
            }
        }
    }
}
